# Applies the changes described in the diff:
# 1. Delete the row for account 004466221 (WALTER, 206168.67)
# 2. Change the balance for account 004352384 (BRASFORT) from 110388.36 to 91111.7
# 3. Delete the row for account 008032413 (VICTOR, 24000)
# 4. Delete the row for account 004395314 (MARIA, 526.19)
# 5. Insert a new row right after account 004451978 (ANTONIO) for
#    account 008002502 (JORGEANA, 500)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Delete WALTER's row (004466221)
$walter = $ws.Columns.Item(1).Find("004466221", $null, $null, 1)
if ($walter -ne $null) {
    $walter.EntireRow.Delete()
}

# 2. Update BRASFORT's balance (004352384): 110388.36 -> 91111.7
$brasfort = $ws.Columns.Item(1).Find("004352384", $null, $null, 1)
if ($brasfort -ne $null) {
    $ws.Cells.Item($brasfort.Row, 3).Value = 91111.7
}

# 3. Delete VICTOR's row (008032413)
$victor = $ws.Columns.Item(1).Find("008032413", $null, $null, 1)
if ($victor -ne $null) {
    $victor.EntireRow.Delete()
}

# 4. Delete MARIA's row (004395314)
$maria = $ws.Columns.Item(1).Find("004395314", $null, $null, 1)
if ($maria -ne $null) {
    $maria.EntireRow.Delete()
}

# 5. Insert JORGEANA's new row (008002502, 500) right after ANTONIO (004451978)
$antonio = $ws.Columns.Item(1).Find("004451978", $null, $null, 1)
$insertRow = $antonio.Row + 1
$ws.Rows.Item($insertRow).Insert()
$ws.Cells.Item($insertRow, 1).Value = "'008002502"
$ws.Cells.Item($insertRow, 2).Value = "JORGEANA"
$ws.Cells.Item($insertRow, 3).Value = 500
